$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# Update the Year value in A2 from 2023 to 2022
$ws.Range("A2").Value = 2022

# Move the active selection to H10
$ws.Activate()
$ws.Range("H10").Select()
